{"js": "// Update the date and the twenty-five three-digit \u00f7 one-digit division\n// problems to the next day's worksheet values.\nconst replacements = [\n  [\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"],\n  [\"716\u00f72=\", \"879\u00f79=\"],\n  [\"284\u00f79=\", \"277\u00f73=\"],\n  [\"844\u00f77=\", \"836\u00f76=\"],\n  [\"642\u00f75=\", \"379\u00f72=\"],\n  [\"887\u00f76=\", \"136\u00f76=\"],\n  [\"888\u00f72=\", \"936\u00f79=\"],\n  [\"409\u00f73=\", \"589\u00f74=\"],\n  [\"961\u00f77=\", \"983\u00f73=\"],\n  [\"733\u00f75=\", \"988\u00f79=\"],\n  [\"588\u00f72=\", \"595\u00f75=\"],\n  [\"306\u00f79=\", \"229\u00f72=\"],\n  [\"358\u00f79=\", \"525\u00f79=\"],\n  [\"181\u00f78=\", \"442\u00f72=\"],\n  [\"607\u00f78=\", \"238\u00f77=\"],\n  [\"489\u00f72=\", \"869\u00f73=\"],\n  [\"787\u00f76=\", \"570\u00f77=\"],\n  [\"648\u00f75=\", \"241\u00f78=\"],\n  [\"251\u00f74=\", \"276\u00f78=\"],\n  [\"547\u00f72=\", \"530\u00f73=\"],\n  [\"410\u00f73=\", \"745\u00f74=\"],\n  [\"412\u00f75=\", \"295\u00f79=\"],\n  [\"398\u00f77=\", \"861\u00f75=\"],\n  [\"231\u00f73=\", \"842\u00f76=\"],\n  [\"995\u00f75=\", \"974\u00f73=\"],\n  [\"755\u00f74=\", \"809\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the twenty-five three-digit \u00f7 one-digit division\n# problems to the next day's worksheet values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"),\n    @(\"716\u00f72=\", \"879\u00f79=\"),\n    @(\"284\u00f79=\", \"277\u00f73=\"),\n    @(\"844\u00f77=\", \"836\u00f76=\"),\n    @(\"642\u00f75=\", \"379\u00f72=\"),\n    @(\"887\u00f76=\", \"136\u00f76=\"),\n    @(\"888\u00f72=\", \"936\u00f79=\"),\n    @(\"409\u00f73=\", \"589\u00f74=\"),\n    @(\"961\u00f77=\", \"983\u00f73=\"),\n    @(\"733\u00f75=\", \"988\u00f79=\"),\n    @(\"588\u00f72=\", \"595\u00f75=\"),\n    @(\"306\u00f79=\", \"229\u00f72=\"),\n    @(\"358\u00f79=\", \"525\u00f79=\"),\n    @(\"181\u00f78=\", \"442\u00f72=\"),\n    @(\"607\u00f78=\", \"238\u00f77=\"),\n    @(\"489\u00f72=\", \"869\u00f73=\"),\n    @(\"787\u00f76=\", \"570\u00f77=\"),\n    @(\"648\u00f75=\", \"241\u00f78=\"),\n    @(\"251\u00f74=\", \"276\u00f78=\"),\n    @(\"547\u00f72=\", \"530\u00f73=\"),\n    @(\"410\u00f73=\", \"745\u00f74=\"),\n    @(\"412\u00f75=\", \"295\u00f79=\"),\n    @(\"398\u00f77=\", \"861\u00f75=\"),\n    @(\"231\u00f73=\", \"842\u00f76=\"),\n    @(\"995\u00f75=\", \"974\u00f73=\"),\n    @(\"755\u00f74=\", \"809\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
